$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 23,7
$arr[0,0] = "Waiver"
$arr[0,1] = "MD"
$arr[0,2] = "1915(c)"
$arr[0,3] = "Amendment"
$arr[0,4] = "MD-2260.R00.55"
$arr[0,5] = "Submitted"
$arr[0,6] = "MD-2260.R00.00"
$arr[1,0] = "SPA"
$arr[1,1] = "MD"
$arr[1,2] = "Medicaid SPA"
$arr[1,3] = ""
$arr[1,4] = "MD-25-9518"
$arr[1,5] = "Approved"
$arr[1,6] = ""
$arr[2,0] = "SPA"
$arr[2,1] = "MD"
$arr[2,2] = "CHIP SPA"
$arr[2,3] = ""
$arr[2,4] = "MD-25-9519"
$arr[2,5] = "Submitted"
$arr[2,6] = ""
$arr[3,0] = "SPA"
$arr[3,1] = "MD"
$arr[3,2] = "Medicaid SPA"
$arr[3,3] = ""
$arr[3,4] = "MD-25-9520"
$arr[3,5] = "Under Review"
$arr[3,6] = ""
$arr[4,0] = "SPA"
$arr[4,1] = "MD"
$arr[4,2] = "Medicaid SPA"
$arr[4,3] = ""
$arr[4,4] = "MD-25-9521"
$arr[4,5] = "Disapproved"
$arr[4,6] = ""
$arr[5,0] = "Waiver"
$arr[5,1] = "MD"
$arr[5,2] = "1915(c)"
$arr[5,3] = "Amendment"
$arr[5,4] = "MD-2260.R00.56"
$arr[5,5] = ""
$arr[5,6] = "MD-2260.R00.00"
$arr[6,0] = "SPA"
$arr[6,1] = "MD"
$arr[6,2] = "Medicaid SPA"
$arr[6,3] = ""
$arr[6,4] = "MD-25-9522"
$arr[6,5] = "Pending-Concurrence"
$arr[6,6] = ""
$arr[7,0] = "SPA"
$arr[7,1] = "MD"
$arr[7,2] = "CHIP SPA"
$arr[7,3] = ""
$arr[7,4] = "MD-25-9523"
$arr[7,5] = "Submitted"
$arr[7,6] = ""
$arr[8,0] = "SPA"
$arr[8,1] = "MD"
$arr[8,2] = "Medicaid SPA"
$arr[8,3] = ""
$arr[8,4] = "MD-25-9524"
$arr[8,5] = "RAI Issued"
$arr[8,6] = ""
$arr[9,0] = "SPA"
$arr[9,1] = "MD"
$arr[9,2] = "CHIP SPA"
$arr[9,3] = ""
$arr[9,4] = "MD-25-9525"
$arr[9,5] = "Submitted"
$arr[9,6] = ""
$arr[10,0] = "SPA"
$arr[10,1] = "MD"
$arr[10,2] = "Medicaid SPA"
$arr[10,3] = ""
$arr[10,4] = "MD-25-9526"
$arr[10,5] = "Submitted"
$arr[10,6] = ""
$arr[11,0] = "SPA"
$arr[11,1] = "MD"
$arr[11,2] = "Medicaid SPA"
$arr[11,3] = ""
$arr[11,4] = "MD-25-9527"
$arr[11,5] = "Submitted"
$arr[11,6] = ""
$arr[12,0] = "Waiver"
$arr[12,1] = "MD"
$arr[12,2] = "1915(c)"
$arr[12,3] = "Amendment"
$arr[12,4] = "MD-2260.R00.57"
$arr[12,5] = "Submitted"
$arr[12,6] = "MD-2260.R00.00"
$arr[13,0] = "SPA"
$arr[13,1] = "MD"
$arr[13,2] = "Medicaid SPA"
$arr[13,3] = ""
$arr[13,4] = "MD-25-9528"
$arr[13,5] = "Under Review"
$arr[13,6] = ""
$arr[14,0] = "Waiver"
$arr[14,1] = "MD"
$arr[14,2] = "1915(b)"
$arr[14,3] = "Initial"
$arr[14,4] = "MD-2281.R00.00"
$arr[14,5] = "Submitted"
$arr[14,6] = ""
$arr[15,0] = "SPA"
$arr[15,1] = "MD"
$arr[15,2] = "Medicaid SPA"
$arr[15,3] = ""
$arr[15,4] = "MD-25-9529"
$arr[15,5] = "Withdrawn"
$arr[15,6] = ""
$arr[16,0] = "SPA"
$arr[16,1] = "MD"
$arr[16,2] = "CHIP SPA"
$arr[16,3] = ""
$arr[16,4] = "MD-25-9530"
$arr[16,5] = "Submitted"
$arr[16,6] = ""
$arr[17,0] = "Waiver"
$arr[17,1] = "MD"
$arr[17,2] = "1915(c)"
$arr[17,3] = "Amendment"
$arr[17,4] = "MD-2260.R00.58"
$arr[17,5] = "Submitted"
$arr[17,6] = "MD-2260.R00.00"
$arr[18,0] = "Waiver"
$arr[18,1] = "MD"
$arr[18,2] = "1915(c)"
$arr[18,3] = "Amendment"
$arr[18,4] = "MD-2260.R00.59"
$arr[18,5] = "Pending-Approval"
$arr[18,6] = "MD-2260.R00.00"
$arr[19,0] = "Waiver"
$arr[19,1] = "MD"
$arr[19,2] = "1915(c)"
$arr[19,3] = "Amendment"
$arr[19,4] = "MD-2260.R00.60"
$arr[19,5] = "Pending-Approval"
$arr[19,6] = "MD-2260.R00.00"
$arr[20,0] = "Waiver"
$arr[20,1] = "MD"
$arr[20,2] = "1915(c)"
$arr[20,3] = "Amendment"
$arr[20,4] = "MD-2260.R00.61"
$arr[20,5] = "Unsubmitted"
$arr[20,6] = "MD-2260.R00.00"
$arr[21,0] = "Waiver"
$arr[21,1] = "MD"
$arr[21,2] = "1915(b)"
$arr[21,3] = "Initial"
$arr[21,4] = "MD-2282.R00.00"
$arr[21,5] = "Terminated"
$arr[21,6] = ""
$arr[22,0] = "Waiver"
$arr[22,1] = "MD"
$arr[22,2] = "1915(b)"
$arr[22,3] = "Initial"
$arr[22,4] = "MD-2283.R00.00"
$arr[22,5] = "Terminated"
$arr[22,6] = ""

$ws.Range("A23:G45").Value = $arr
